$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-empty "Address" cells in column F (they contained no text)
$ws.Range("F3").ClearContents()
$ws.Range("F22").ClearContents()
$ws.Range("F42").ClearContents()
$ws.Range("F49").ClearContents()

# Update District names (column G) to their official names:
# "Tumkur" and "Madhugiri" both become "Tumakuru (Tumkur)"
for ($r = 4; $r -le 57; $r++) {
    $cell = $ws.Range("G$r")
    $val = $cell.Value2
    if ($val -eq "Tumkur" -or $val -eq "Madhugiri") {
        $cell.Value = "Tumakuru (Tumkur)"
    }
}
